$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Two outlier measurements (Picture numbers 19 and 42) were dropped from the
# first results table. Every later row of that table shifts up by one row
# for each removed entry, while the second table (starting at row 26) must
# stay exactly where it is. So instead of deleting whole sheet rows (which
# would also push the second table upward), we overwrite the affected cells
# directly and clear out what is left over at the bottom of the first table.

# Rows 3-5 (Picture 15, 17, 18) are unchanged.

# Row 6 <- old row 7 (Picture 20)
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = 20
$ws.Cells.Item(6, 3).Value = 0.97936640865681601
$ws.Cells.Item(6, 4).Value = 3.7735924528225802
$ws.Cells.Item(6, 5).Value = 0.14646899999999999

# Row 7 <- old row 8 (Picture 24)
$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = 24
$ws.Cells.Item(7, 3).Value = 0.61367503537374601
$ws.Cells.Item(7, 4).Value = 8.0622577482985491
$ws.Cells.Item(7, 5).Value = 0.14398939999999999

# Row 8 <- old row 9 (Picture 32)
$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = 32
$ws.Cells.Item(8, 3).Value = 0.26248787031142801
$ws.Cells.Item(8, 4).Value = 0.76987447698745404
$ws.Cells.Item(8, 5).Value = 0.1797899

# Row 9 <- old row 10 (Picture 33)
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = 33
$ws.Cells.Item(9, 3).Value = 0.17433204441147199
$ws.Cells.Item(9, 4).Value = 0.294520547945197
$ws.Cells.Item(9, 5).Value = 0.12534300000000001

# Row 10 <- old row 11 (Picture 34)
$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = 34
$ws.Cells.Item(10, 3).Value = 0.097711744135551698
$ws.Cells.Item(10, 4).Value = 0.57575757575757802
$ws.Cells.Item(10, 5).Value = 0.1378943

# Row 11 <- old row 12 (Picture 38)
$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = 38
$ws.Cells.Item(11, 3).Value = 0.15869675698938501
$ws.Cells.Item(11, 4).Value = 0.064814814814809593
$ws.Cells.Item(11, 5).Value = 0.1391773

# Row 12 <- old row 13 (Picture 39)
$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = 39
$ws.Cells.Item(12, 3).Value = 1.53979563663623
$ws.Cells.Item(12, 4).Value = 4.1481481481481497
$ws.Cells.Item(12, 5).Value = 0.17536740000000001

# Row 13 <- old row 14 (Picture 40)
$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = 40
$ws.Cells.Item(13, 3).Value = 0.15156997140904099
$ws.Cells.Item(13, 4).Value = 0.093457943925230297
$ws.Cells.Item(13, 5).Value = 0.16435449999999999

# Row 14 <- old row 15 (Picture 41)
$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = 41
$ws.Cells.Item(14, 3).Value = 0.151648285813309
$ws.Cells.Item(14, 4).Value = 0.026119402985074199
$ws.Cells.Item(14, 5).Value = 0.1797559

# Row 15 <- old row 17 (Picture 43) -- old row 16 (Picture 42) dropped
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = 43
$ws.Cells.Item(15, 3).Value = 0.69545276946352896
$ws.Cells.Item(15, 4).Value = 0.954545454545496
$ws.Cells.Item(15, 5).Value = 0.17295559999999999

# Row 16 <- old row 18 (Picture 44)
$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = 44
$ws.Cells.Item(16, 3).Value = 0.16248989705060499
$ws.Cells.Item(16, 4).Value = 2.0750000000000499
$ws.Cells.Item(16, 5).Value = 0.14608959999999999

# Row 17 <- old row 19 (Picture 45)
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = 45
$ws.Cells.Item(17, 3).Value = 0.26122514825186699
$ws.Cells.Item(17, 4).Value = 1.0657082490166601
$ws.Cells.Item(17, 5).Value = 0.14041790000000001

# Old data rows 18 and 19 no longer hold data now that the table ends at
# row 17 - clear them out.
$ws.Range("A18:E19").ClearContents()

# The min/max/avg summary rows move up from 21-23 to 19-21, with formula
# ranges now covering only the remaining 15 data rows (3:17).
$ws.Cells.Item(19, 2).Value = "min"
$ws.Cells.Item(19, 3).Formula = "=MIN(C3:C17)"
$ws.Cells.Item(19, 4).Formula = "=MIN(D3:D17)"
$ws.Cells.Item(19, 5).Formula = "=MIN(E3:E17)"

$ws.Cells.Item(20, 2).Value = "max"
$ws.Cells.Item(20, 3).Formula = "=MAX(C3:C17)"
$ws.Cells.Item(20, 4).Formula = "=MAX(D3:D17)"
$ws.Cells.Item(20, 5).Formula = "=MAX(E3:E17)"

$ws.Cells.Item(21, 2).Value = "avg"
$ws.Cells.Item(21, 3).Formula = "=AVERAGE(C3:C17)"
$ws.Cells.Item(21, 4).Formula = "=AVERAGE(D3:D17)"
$ws.Cells.Item(21, 5).Formula = "=AVERAGE(E3:E17)"

# The old summary rows (22 and 23) are no longer used now that the
# summary moved up to rows 19-21; clear them so the gap before the
# second table (row 26) matches.
$ws.Range("A22:E23").ClearContents()

$ws.Range("J14").Select()
